$d = $word.ActiveDocument

# The document's headers contain the BTec orange logo picture (currently
# named "image1.jpg") which must be renamed to "image2.jpg", and the
# footers contain the Pearson logo picture (currently named "image2.png")
# which must be renamed to "image1.png". Walk every section's headers and
# footers (primary, first-page, even-page) and every inline picture in
# them, and rename based on which logo it is (identified by its
# AlternativeText / description, which is preserved by the edit).

for ($sIdx = 1; $sIdx -le $d.Sections.Count; $sIdx++) {
    $sec = $d.Sections.Item($sIdx)

    for ($h = 1; $h -le 3; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    for ($f = 1; $f -le 3; $f++) {
        $ftr = $sec.Footers.Item($f)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
